$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$newStatus = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- 2. Latest Handback DateTime placeholder -> real timestamp (zh-cn) ---
$zhcn.Range("K2").Value = "2016-09-03 06:31:57"
$zhcn.Range("K3").Value = "2016-09-03 06:31:57"

# --- 3. zh-cn: fill Latest Target File (I) and Latest Handback File (J) ---
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c4692c86d76666e880cdaee664749078b52d5050/e2e/8dda6b2d-950c-4290-804b-b6de9f7e2641.md", "", "", "8dda6b2d-950c-4290-804b-b6de9f7e2641.md")
$zhcn.Range("J2").Value = "8dda6b2d-950c-4290-804b-b6de9f7e2641.05cd43248fd0e34d34b7e27d332ee368e4b09646.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c4692c86d76666e880cdaee664749078b52d5050/e2e/b4ee3680-7ceb-47c5-aee1-656d7f283a8f.md", "", "", "b4ee3680-7ceb-47c5-aee1-656d7f283a8f.md")
$zhcn.Range("J3").Value = "b4ee3680-7ceb-47c5-aee1-656d7f283a8f.fd165c579508646209e570f3ec1a3faeec0d3d83.zh-cn.xlf"

# --- 4. de-de: fill Latest Handoff File (G), Latest Target File (I), Latest Handback File (J),
#        Latest Handback DateTime (K) ---
$dede.Range("G2").Value = "2016-09-03 06:32:09"
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c4692c86d76666e880cdaee664749078b52d5050/e2e/8dda6b2d-950c-4290-804b-b6de9f7e2641.md", "", "", "8dda6b2d-950c-4290-804b-b6de9f7e2641.md")
$dede.Range("J2").Value = "2016-09-03 06:32:09"
$dede.Range("K2").Value = "8dda6b2d-950c-4290-804b-b6de9f7e2641.05cd43248fd0e34d34b7e27d332ee368e4b09646.de-de.xlf"

$dede.Range("G3").Value = "b4ee3680-7ceb-47c5-aee1-656d7f283a8f.fd165c579508646209e570f3ec1a3faeec0d3d83.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c4692c86d76666e880cdaee664749078b52d5050/e2e/b4ee3680-7ceb-47c5-aee1-656d7f283a8f.md", "", "", "b4ee3680-7ceb-47c5-aee1-656d7f283a8f.md")
$dede.Range("J3").Value = "b4ee3680-7ceb-47c5-aee1-656d7f283a8f.fd165c579508646209e570f3ec1a3faeec0d3d83.de-de.xlf"
$dede.Range("K3").Value = "8dda6b2d-950c-4290-804b-b6de9f7e2641.05cd43248fd0e34d34b7e27d332ee368e4b09646.de-de.xlf"

# --- 5. Column width adjustments ---
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40
